$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.145.22"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "2.525.91"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.94%  "

$ws.Range("D9").Value = "2.525.91"
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("E11").Value = "  -2.11%  "

$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").Value = "2.970.93"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "59.059.40"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "2.527.85"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("E26").Value = "  -1.80%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("D30").Value = "0.0₃0774"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "

$ws.Range("E33").Value = "  +5.60%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "

$ws.Range("E40").Value = "  +0.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "284.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.70%  "

$ws.Range("E46").Value = "  +2.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0924"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("E50").Value = "  -1.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
